$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style of the existing header row (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the new numeric columns I and J for rows 2-18
$data = @(
    @(9, 9),
    @(6, 8),
    @(7, 8),
    @(9, 9),
    @(4, 7),
    @(8, 9),
    @(8, 8),
    @(7, 8),
    @(7, 7),
    @(7, 7),
    @(4, 6),
    @(4, 8),
    @(9, 9),
    @(7, 9),
    @(4, 6),
    @(5, 6),
    @(3, 4)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
